# Fruta / hortaliza, semanal
#
# The edit reshuffles the per-record data (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Origen, Precio $/Kg) across the data
# rows (2-43): every row keeps its "constant" columns (Mercado ID,
# Mercado, Region, Codreg, Categoria ID, Categoria, Variedad, Unidad de
# comercializacion, Kg o Unidades, Clasificacion) but the "variable"
# columns are redistributed to different rows. Because this is a
# genuine row-to-row permutation, we must snapshot every source row's
# values BEFORE writing anything, then apply the new layout from that
# snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move between rows.
$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

$firstRow = 2
$lastRow = 43

# 1) Snapshot the current ("before") values for every data row.
$orig = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $row
}

# 2) Destination row -> source row mapping (the permutation taken from
#    the target workbook).
$mapping = @{
    2  = 38
    3  = 27
    4  = 18
    5  = 29
    6  = 9
    7  = 17
    8  = 34
    9  = 25
    10 = 24
    11 = 37
    12 = 36
    13 = 14
    14 = 15
    15 = 31
    16 = 30
    17 = 41
    18 = 11
    19 = 12
    20 = 3
    21 = 43
    22 = 16
    23 = 13
    24 = 8
    25 = 10
    26 = 19
    27 = 7
    28 = 22
    29 = 39
    30 = 32
    31 = 6
    32 = 40
    33 = 28
    34 = 33
    35 = 4
    36 = 42
    37 = 5
    38 = 21
    39 = 26
    40 = 23
    41 = 20
    42 = 2
    43 = 35
}

# 3) Write the permuted values back out using the snapshot, so the
#    source/destination overlap in row numbers does not clobber data
#    still waiting to be read.
foreach ($dest in $mapping.Keys) {
    $srcRow = $mapping[$dest]
    $srcData = $orig[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value = $srcData[$c]
    }
}
